# Insert a new data row at row 62 (pushing existing rows 62-92 down to 63-93)
# and populate it with a new "Poroto verde" price-report record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(62).Insert()

$ws.Cells.Item(62, 1).Value = 10
$ws.Cells.Item(62, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(62, 3).Value = "La Araucanía"
$ws.Cells.Item(62, 4).Value = 44523
$ws.Cells.Item(62, 5).Value = 9
$ws.Cells.Item(62, 6).Value = 100112031
$ws.Cells.Item(62, 7).Value = "Poroto verde"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 125
$ws.Cells.Item(62, 11).Value = 2000
$ws.Cells.Item(62, 12).Value = 2000
$ws.Cells.Item(62, 13).Value = 2000
$ws.Cells.Item(62, 14).Value = "$/kilo"
$ws.Cells.Item(62, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(62, 16).Value = 2000
$ws.Cells.Item(62, 17).Value = 1
$ws.Cells.Item(62, 18).Value = "Hortaliza"
